$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Índice"
$ws.Range("B1").Value = "Distancia"
$ws.Range("C1").Value = "max"
$ws.Range("D1").Value = "min"
$ws.Range("E1").Value = "Tempo"

$data = @(
    @(0, 3435.566666666667, 3512, 3399, 0.02883265018463135),
    @(1, 3234.333333333333, 3387, 3029, 0.03097601731618245),
    @(2, 3458.2, 3585, 3238, 0.03085882663726807),
    @(3, 3570.733333333333, 3917, 3312, 0.03063379128774007),
    @(4, 3298, 3555, 2981, 0.03138500849405924),
    @(5, 4081.466666666667, 4227, 3990, 0.03108121554056803),
    @(6, 3107.2, 3344, 2966, 0.03072756131490072),
    @(7, 3790.366666666667, 4021, 3582, 0.03126941521962483),
    @(8, 3149.266666666667, 3377, 2935, 0.03217519124348958),
    @(9, 3122.9, 3531, 3029, 0.02890580495198568)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
    $ws.Cells.Item($row, 4).Value = $data[$i][3]
    $ws.Cells.Item($row, 5).Value = $data[$i][4]
}
